$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some target values are numeric-looking strings (e.g. "1.001") but must
# remain stored as text, matching the source data. Force text format on
# just those cells before assigning so Excel does not auto-convert them
# to numbers.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.221.93"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "1.907.18"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "308.27"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D7").Value = "0.5244"
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").Value = "0.3781"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("D9").Value = "0.07274"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "21.27"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").Value = "0.8994"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "0.07683"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "1.918.76"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "95.27"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "5.273"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "0.000008667"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "14.54"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "27.298.06"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").Value = "2.148.66"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("D24").Value = "6.449"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "2.326"
$ws.Range("E25").Value = "  +11.32%  "
$ws.Range("D26").Value = "145.59"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "18.15"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "1.737"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").Value = "114.89"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").Value = "4.971"
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("D31").Value = "4.817"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").Value = "0.09239"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.8088"
$ws.Range("E33").Value = "  +8.07%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.05077"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "1.245"
$ws.Range("E35").Value = "  +7.90%  "
$ws.Range("D36").Value = "2.996"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "3.315"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").Value = "2.603"
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("D39").Value = "0.5674"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "1.076"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "9.006"
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("D43").Value = "6.644"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "119.31"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "0.1519"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").Value = "0.4846"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "10.27"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "1.618"
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("D50").Value = "37.55"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "63.95"
$ws.Range("E51").Value = "  +1.40%  "
